$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add a new row (14) documenting the "Oil Prices" time-series forecasting
# project (Stacked Ensemble: Random Forest Regressor meta-learner on SARIMAX
# and Holt-Winters models).
# ---------------------------------------------------------------------------

# --- Formatting -------------------------------------------------------
# A14 needs the same look as A11 (italic Palatino, left/center, wrap text)
# but with a full thin box border (A11's border is missing its left edge
# because it anchors a vertical merge; A14 is a standalone row).
$ws.Range("A11").Copy() | Out-Null
$ws.Range("A14").PasteSpecial(-4122) | Out-Null
$ws.Range("A14").Borders.Item(7).LineStyle = 1
$ws.Range("A14").Borders.Item(7).Weight = 2

# B14 mirrors the other "Type of ML Model" cells (e.g. B7 / B11 / B12).
$ws.Range("B11").Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4122) | Out-Null

# C14 mirrors the other hyperlink "Project Name" cells (e.g. C12 / C13).
$ws.Range("C12").Copy() | Out-Null
$ws.Range("C14").PasteSpecial(-4122) | Out-Null

$ws.Application.CutCopyMode = 0

# --- Values -------------------------------------------------------------
$ws.Range("A14").Value = "Stacked Ensemble (Random Forest Regressor Meta-Learner on SARIMAX and Holt-Winters Models)."
$ws.Range("B14").Value = "Time Series Forecasting"
$ws.Range("C14").Value = "Oil Prices"

# --- Row height (matches A11's two-line wrapped height) -----------------
$ws.Rows(14).RowHeight = 31.2

# --- Hyperlink ------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("C14"), "c. Jupyter Notebooks\Oil Prices.ipynb", "", "", "Oil Prices") | Out-Null

# --- Selection (matches the saved view state in the target workbook) ------
$ws.Range("A6").Select() | Out-Null
